$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few more practice programs got finished -- fill in the file name (col B)
# and mark them done (col C) for the matching rows in col A.
$ws.Range("B8").Value  = "swapWithout3Variable.java"
$ws.Range("C8").Value  = "DONE"

$ws.Range("B9").Value  = "evenOrOdd.java"
$ws.Range("C9").Value  = "DONE"

$ws.Range("B10").Value = "vowelOrConsonant.java"
$ws.Range("C10").Value = "DONE"

$ws.Range("B14").Value = "positiveOrNegative.java"
$ws.Range("C14").Value = "DONE"

$ws.Range("B15").Value = "aplhabetOrNot.java"
$ws.Range("C15").Value = "DONE"

$ws.Range("B16").Value = "sumOfNaturalNumbers.java"
$ws.Range("C16").Value = "DONE"

$ws.Range("B11").Value = "largestAmongThree.java"
$ws.Range("C11").Value = "DONE"

# Leave the view where the author left it: scrolled down a bit further,
# with C11 as the active/selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

$ws.Range("C11").Select()
